$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each touched cell is forced to Text format before the value is written,
# so numeric-looking strings (e.g. prices) are preserved exactly as text
# instead of being auto-converted to floating point numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.716.29'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.446.17'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.95%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.06'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.58'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.62%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.445.42'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.64%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.78%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.00%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.19'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +7.84%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.568.83'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.16%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.448.26'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.84'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.83%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.94'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.55%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'BabyDogeCoin'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₆0836'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +198.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '329.97'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.86%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +9.44%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '642.68'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +12.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.17'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +17.61%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +5.23%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.565.80'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.46'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +8.74%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.51'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +6.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.374'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '152.24'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '18.71'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.70'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +5.70%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '42.23'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '14.94'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +27.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '145.23'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.04%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.62'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +6.55%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.98%  '
